# Add an "Acknowledgments" heading (with bookmark) before the Credits
# paragraph, update the Credits paragraph text, and give the (previously
# empty) final section explicit page setup (matches Word defaults).

$d = $word.ActiveDocument

# Locate the "Credits" paragraph (currently holds the placeholder text).
$targetIndex = 0
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Style.NameLocal -eq "Credits") {
        $targetIndex = $idx
    }
}

$creditsPara = $d.Paragraphs.Item($targetIndex)

# Insert a new paragraph right before it for the "Acknowledgments" heading.
$creditsPara.Range.InsertParagraphBefore()

$headingPara = $d.Paragraphs.Item($targetIndex)
$headingPara.Style = "Heading1"
$headingPara.Range.Text = "Acknowledgments"

# Wrap the heading text in a bookmark named "acknowledgments".
$bmRange = $d.Range($headingPara.Range.Start, $headingPara.Range.Start + 15)
$d.Bookmarks.Add("acknowledgments", $bmRange)

# Update the Credits paragraph text (it shifted down by one paragraph).
$creditsPara = $d.Paragraphs.Item($targetIndex + 1)
$creditsPara.Range.Text = "Some materials included in this export came from the following casebooks."

# Give the last section an explicit page setup (Word-default Letter page,
# 1in margins, 0.5in header/footer distance, no gutter, default column
# spacing) instead of the empty <w:sectPr/>.
$ps = $d.Sections.Item($d.Sections.Count).PageSetup
$ps.PageWidth = 612
$ps.PageHeight = 792
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.LeftMargin = 72
$ps.RightMargin = 72
$ps.HeaderDistance = 36
$ps.FooterDistance = 36
$ps.Gutter = 0
$ps.TextColumns.Spacing = 36
